# Apply updated dSF (column F) values as per repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = 3
    4  = -2
    5  = 8
    6  = -1
    7  = 1
    8  = 4
    10 = 3
    11 = 4
    14 = 5
    15 = 2
    16 = -4
    18 = -3
    19 = 2
    20 = 3
    21 = 4
    23 = -1
    24 = 1
    26 = -2
    27 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
